# Fruta / hortaliza, semanal
# Re-point the weekly rows (2-11) at the correct Fecha / Volumen / Precio
# values - the underlying facts were shuffled across rows and this
# restores the correct row-to-data alignment (row 4 is already correct).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @{ D = 44761; M = 100; N = 20000; O = 21000; P = 20500; S = 1025 }
    3  = @{ D = 44357; M = 100; N = 14000; O = 15000; P = 14500; S = 725 }
    4  = @{ D = 44798; M = 80;  N = 21000; O = 22000; P = 21500; S = 1075 }
    5  = @{ D = 44708; M = 80;  N = 20000; O = 21000; P = 20500; S = 1025 }
    6  = @{ D = 44533; M = 100; N = 16000; O = 17000; P = 16500; S = 825 }
    7  = @{ D = 45092; M = 150; N = 24000; O = 25000; P = 24333; S = 1217 }
    8  = @{ D = 44320; M = 80;  N = 16000; O = 17000; P = 16500; S = 825 }
    9  = @{ D = 44792; M = 100; N = 21000; O = 22000; P = 21500; S = 1075 }
    10 = @{ D = 44893; M = 80;  N = 21000; O = 22000; P = 21625; S = 1081 }
    11 = @{ D = 44890; M = 80;  N = 20000; O = 23000; P = 22250; S = 1112 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("M$row").Value = $vals.M
    $ws.Range("N$row").Value = $vals.N
    $ws.Range("O$row").Value = $vals.O
    $ws.Range("P$row").Value = $vals.P
    $ws.Range("S$row").Value = $vals.S
}
